# Auto commit at 2025-12-03  7:39:20.31
#
# Refresh the "Metrics" sheet with the latest source numbers, then re-select
# the working cells on both the Metrics and today sheets (the today sheet's
# TODAY()-1 header cell and its lookups recalc automatically).

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

# Each metric grew by the same per-row delta it had in the prior refresh
# (kwh/service-fee deltas repeat across the monthly/yearly roll-up rows),
# so add the delta to the existing cached value instead of retyping a
# rounded literal - this keeps the float bit pattern identical to a real
# recalculated roll-up.
$metrics.Range("B2").Value = $metrics.Range("B2").Value() + 13589.65
$metrics.Range("B3").Value = $metrics.Range("B3").Value() + 11786.04
$metrics.Range("B4").Value = $metrics.Range("B4").Value() + 4089.23
$metrics.Range("B5").Value = $metrics.Range("B5").Value() + 563
$metrics.Range("B6").Value = $metrics.Range("B6").Value() + 13589.65
$metrics.Range("B7").Value = $metrics.Range("B7").Value() + 11786.04
$metrics.Range("B8").Value = $metrics.Range("B8").Value() + 4089.23
$metrics.Range("B9").Value = $metrics.Range("B9").Value() + 563
$metrics.Range("B10").Value = $metrics.Range("B10").Value() + 13589.65
$metrics.Range("B11").Value = $metrics.Range("B11").Value() + 11786.04
$metrics.Range("B12").Value = $metrics.Range("B12").Value() + 4089.23
$metrics.Range("B13").Value = $metrics.Range("B13").Value() + 563

$metrics.Activate()
$metrics.Range("D7").Select()

$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("F9").Select()
